$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("L2").Value = 0.93
$ws1.Range("L3").Value = 0.95
$ws1.Range("L4").Value = 1.18

$ws1.Range("D5").Value = 1
$ws1.Range("H5").Value = 0.5
$ws1.Range("I5").Value = "High"
$ws1.Range("L5").Value = 1.2

$ws1.Range("L6").Value = 1.14
$ws1.Range("L7").Value = 1.01
$ws1.Range("L8").Value = 1.18
$ws1.Range("L9").Value = 0.85

$ws1.Range("L11").Value = 1.05
$ws1.Range("L12").Value = 0.93
$ws1.Range("L13").Value = 0.99
$ws1.Range("L14").Value = 0.95
$ws1.Range("L15").Value = 0.98
$ws1.Range("L16").Value = 0.95
$ws1.Range("L17").Value = 0.85

# --- Sheet "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B10").Value = "'8"
$ws2.Range("B11").Value = "'4"
